$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, [string]$Text)
    $Cell.Value = "'" + $Text
    $Cell.Style = "Normal"
}

Set-TextValue $ws.Cells.Item(2, 4) "43.875.66"
Set-TextValue $ws.Cells.Item(2, 5) "  +0.60%  "
Set-TextValue $ws.Cells.Item(3, 4) "2.342.54"
Set-TextValue $ws.Cells.Item(3, 5) "  +4.74%  "
Set-TextValue $ws.Cells.Item(4, 5) "  +0.04%  "
Set-TextValue $ws.Cells.Item(5, 4) "96.99"
Set-TextValue $ws.Cells.Item(5, 5) "  +2.88%  "
Set-TextValue $ws.Cells.Item(6, 4) "272.10"
Set-TextValue $ws.Cells.Item(6, 5) "  +0.92%  "
Set-TextValue $ws.Cells.Item(7, 5) "  +0.55%  "
Set-TextValue $ws.Cells.Item(8, 5) "  -0.03%  "
Set-TextValue $ws.Cells.Item(9, 4) "0.628"
Set-TextValue $ws.Cells.Item(9, 5) "  +0.44%  "
Set-TextValue $ws.Cells.Item(10, 4) "45.78"
Set-TextValue $ws.Cells.Item(10, 5) "  -1.10%  "
Set-TextValue $ws.Cells.Item(11, 4) "0.0948"
Set-TextValue $ws.Cells.Item(11, 5) "  +2.54%  "
Set-TextValue $ws.Cells.Item(12, 4) "8.16"
Set-TextValue $ws.Cells.Item(12, 5) "  +0.00%  "
Set-TextValue $ws.Cells.Item(13, 5) "  +0.53%  "
Set-TextValue $ws.Cells.Item(14, 4) "2.690.67"
Set-TextValue $ws.Cells.Item(14, 5) "  +4.62%  "
Set-TextValue $ws.Cells.Item(15, 4) "15.68"
Set-TextValue $ws.Cells.Item(15, 5) "  +3.38%  "
Set-TextValue $ws.Cells.Item(16, 4) "0.871"
Set-TextValue $ws.Cells.Item(16, 5) "  +9.03%  "
Set-TextValue $ws.Cells.Item(17, 4) "2.339.77"
Set-TextValue $ws.Cells.Item(17, 5) "  +4.28%  "
Set-TextValue $ws.Cells.Item(18, 4) "43.847.37"
Set-TextValue $ws.Cells.Item(18, 5) "  +0.60%  "
Set-TextValue $ws.Cells.Item(19, 4) "0.0000109"
Set-TextValue $ws.Cells.Item(19, 5) "  +5.37%  "
Set-TextValue $ws.Cells.Item(20, 4) "6.46"
Set-TextValue $ws.Cells.Item(20, 5) "  +7.15%  "
Set-TextValue $ws.Cells.Item(21, 4) "72.91"
Set-TextValue $ws.Cells.Item(21, 5) "  +3.43%  "
Set-TextValue $ws.Cells.Item(22, 4) "240.21"
Set-TextValue $ws.Cells.Item(22, 5) "  +2.88%  "
Set-TextValue $ws.Cells.Item(23, 5) "  -1.97%  "
Set-TextValue $ws.Cells.Item(24, 5) "  +4.91%  "
Set-TextValue $ws.Cells.Item(25, 5) "  -0.10%  "
Set-TextValue $ws.Cells.Item(26, 5) "  +1.27%  "
Set-TextValue $ws.Cells.Item(27, 4) "11.44"
Set-TextValue $ws.Cells.Item(27, 5) "  +1.65%  "
Set-TextValue $ws.Cells.Item(28, 5) "  -1.70%  "
Set-TextValue $ws.Cells.Item(29, 4) "2.27"
Set-TextValue $ws.Cells.Item(29, 5) "  +0.18%  "
Set-TextValue $ws.Cells.Item(30, 4) "22.60"
Set-TextValue $ws.Cells.Item(30, 5) "  +8.51%  "
Set-TextValue $ws.Cells.Item(31, 4) "38.13"
Set-TextValue $ws.Cells.Item(31, 5) "  -5.49%  "
Set-TextValue $ws.Cells.Item(32, 4) "173.88"
Set-TextValue $ws.Cells.Item(32, 5) "  +0.50%  "
Set-TextValue $ws.Cells.Item(33, 5) "  -2.84%  "
Set-TextValue $ws.Cells.Item(34, 4) "5.50"
Set-TextValue $ws.Cells.Item(34, 5) "  +0.58%  "
Set-TextValue $ws.Cells.Item(35, 4) "0.127"
Set-TextValue $ws.Cells.Item(35, 5) "  +2.45%  "
Set-TextValue $ws.Cells.Item(36, 4) "0.0360"
Set-TextValue $ws.Cells.Item(36, 5) "  +2.66%  "
Set-TextValue $ws.Cells.Item(37, 5) "  -2.14%  "
Set-TextValue $ws.Cells.Item(38, 4) "4.42"
Set-TextValue $ws.Cells.Item(38, 5) "  +2.16%  "
Set-TextValue $ws.Cells.Item(39, 4) "3.40"
Set-TextValue $ws.Cells.Item(39, 5) "  -5.38%  "
Set-TextValue $ws.Cells.Item(40, 4) "2.40"
Set-TextValue $ws.Cells.Item(40, 5) "  +10.06%  "
Set-TextValue $ws.Cells.Item(41, 4) "0.240"
Set-TextValue $ws.Cells.Item(41, 5) "  +9.50%  "
Set-TextValue $ws.Cells.Item(42, 4) "1.39"
Set-TextValue $ws.Cells.Item(42, 5) "  +19.86%  "
Set-TextValue $ws.Cells.Item(43, 4) "12.19"
Set-TextValue $ws.Cells.Item(43, 5) "  -3.66%  "
Set-TextValue $ws.Cells.Item(44, 4) "9.28"
Set-TextValue $ws.Cells.Item(44, 5) "  +10.67%  "
Set-TextValue $ws.Cells.Item(45, 4) "62.54"
Set-TextValue $ws.Cells.Item(45, 5) "  -0.99%  "
Set-TextValue $ws.Cells.Item(46, 4) "5.39"
Set-TextValue $ws.Cells.Item(46, 5) "  +0.68%  "
Set-TextValue $ws.Cells.Item(47, 5) "  +5.44%  "
Set-TextValue $ws.Cells.Item(48, 4) "100.91"
Set-TextValue $ws.Cells.Item(48, 5) "  +0.32%  "
Set-TextValue $ws.Cells.Item(49, 5) "  +1.25%  "
Set-TextValue $ws.Cells.Item(50, 4) "2.570.09"
Set-TextValue $ws.Cells.Item(50, 5) "  +4.55%  "
Set-TextValue $ws.Cells.Item(51, 4) "0.186"
Set-TextValue $ws.Cells.Item(51, 5) "  +14.55%  "
